$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1590.8182
$ws.Range("J29").Value = 1899.8889
$ws.Range("L29").Value = 5699.6667
$ws.Range("N29").Value = -6261.6667

$ws.Range("H32").Value = 1360
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents() | Out-Null

$ws.Range("H33").Value = 14707715
$ws.Range("I33").Value = 19231628
$ws.Range("K33").Value = 19231628
$ws.Range("M33").Value = -19231399

$ws.Range("H64").Value = 6948.8335
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 6948.8335
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 6948.8335
$ws.Range("M64").ClearContents() | Out-Null
$ws.Range("N64").Value = -7444.8335

$ws.Range("H67").Value = 6948.8335
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 6948.8335
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 6948.8335
$ws.Range("M67").ClearContents() | Out-Null
$ws.Range("N67").Value = -8664.833500000001

$ws.Range("H98").Value = 4688.304
$ws.Range("I98").Value = 4602.6113
$ws.Range("K98").Value = 4602.6113
$ws.Range("M98").Value = -3104.6113

$ws.Range("H112").Value = 10684.4
$ws.Range("J112").Value = 4093.7778
$ws.Range("L112").Value = 12281.3334
$ws.Range("N112").Value = -14497.3334

$ws.Range("H113").Value = 8987.529
$ws.Range("I113").Value = 11769.7
$ws.Range("J113").Value = 5013
$ws.Range("K113").Value = 11769.7
$ws.Range("L113").Value = 5013
$ws.Range("M113").Value = -8515.700000000001
$ws.Range("N113").Value = -11521

$ws.Range("H122").Value = 4688.304
$ws.Range("I122").Value = 4602.6113
$ws.Range("K122").Value = 13807.8339
$ws.Range("M122").Value = -11357.8339

$ws.Range("H138").Value = 4087.2363
$ws.Range("I138").Value = 4230.5
$ws.Range("J138").Value = 4069.6938
$ws.Range("K138").Value = 12691.5
$ws.Range("L138").Value = 12209.0814
$ws.Range("M138").Value = -7551.5
$ws.Range("N138").Value = -22489.0814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4244.075
$ws.Range("I110").Value = 5037.6333
$ws.Range("K110").Value = 5037.6333
$ws.Range("M110").Value = -2992.6333

$ws.Range("H122").Value = 3852.7646
$ws.Range("I122").Value = 1549.8
$ws.Range("K122").Value = 4649.4
$ws.Range("M122").Value = -2199.4

$ws.Range("H132").Value = 6308929.5
$ws.Range("I132").Value = 1620.625
$ws.Range("J132").Value = 26492318
$ws.Range("K132").Value = 4861.875
$ws.Range("L132").Value = 79476954
$ws.Range("M132").Value = -2331.875
$ws.Range("N132").Value = -79482014

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 14548
$ws.Range("I23").Value = 200
$ws.Range("K23").Value = 200
$ws.Range("M23").Value = 83

$ws.Range("H107").Value = 1425.5385
$ws.Range("I107").Value = 1458.5385
$ws.Range("J107").Value = 1392.5385
$ws.Range("K107").Value = 1458.5385
$ws.Range("L107").Value = 1392.5385
$ws.Range("M107").Value = 461.4614999999999
$ws.Range("N107").Value = -5232.538500000001

$ws.Range("H134").Value = 41652.887
$ws.Range("I134").Value = 58622.42
$ws.Range("K134").Value = 175867.26
$ws.Range("M134").Value = -173332.26

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 501080.5
$ws.Range("J12").Value = 1666931.6
$ws.Range("L12").Value = 1666931.6
$ws.Range("N12").Value = -1667271.6

$ws.Range("H31").Value = 11134.9
$ws.Range("J31").Value = 41548.2
$ws.Range("L31").Value = 41548.2
$ws.Range("N31").Value = -42138.2

$ws.Range("H34").Value = 11134.9
$ws.Range("J34").Value = 41548.2
$ws.Range("L34").Value = 41548.2
$ws.Range("N34").Value = -41952.2

$ws.Range("H86").Value = 10118.117
$ws.Range("I86").Value = 10200.6
$ws.Range("J86").Value = 9499.5
$ws.Range("K86").Value = 10200.6
$ws.Range("L86").Value = 9499.5
$ws.Range("M86").Value = -9077.6
$ws.Range("N86").Value = -11745.5

$ws.Range("H89").Value = 10118.117
$ws.Range("I89").Value = 10200.6
$ws.Range("J89").Value = 9499.5
$ws.Range("K89").Value = 51003
$ws.Range("L89").Value = 47497.5
$ws.Range("M89").Value = -45387
$ws.Range("N89").Value = -58729.5

$ws.Range("H105").Value = 9926.23
$ws.Range("I105").Value = 15290
$ws.Range("J105").Value = 3668.5
$ws.Range("K105").Value = 15290
$ws.Range("L105").Value = 3668.5
$ws.Range("M105").Value = -13543
$ws.Range("N105").Value = -7162.5

$ws.Range("H132").Value = 49022228
$ws.Range("J132").Value = 205885860
$ws.Range("L132").Value = 617657580
$ws.Range("N132").Value = -617662640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3192201
$ws.Range("J4").Value = 3001249.8
$ws.Range("L4").Value = 9003749.399999999
$ws.Range("N4").Value = -9003973.399999999

$ws.Range("H39").Value = 4177.6665
$ws.Range("J39").Value = 4548.9414
$ws.Range("L39").Value = 13646.8242
$ws.Range("N39").Value = -14234.8242

$ws.Range("H55").Value = 1838.4445
$ws.Range("I55").Value = 1849.6666
$ws.Range("J55").Value = 1816
$ws.Range("K55").Value = 5548.9998
$ws.Range("L55").Value = 5448
$ws.Range("M55").Value = -5371.9998
$ws.Range("N55").Value = -5802

$ws.Range("H98").Value = 788.25
$ws.Range("I98").Value = 803
$ws.Range("K98").Value = 2409
$ws.Range("M98").Value = -911

$ws.Range("H131").Value = 1452.27
$ws.Range("I131").Value = 967.125
$ws.Range("J131").Value = 1494.4565
$ws.Range("K131").Value = 2901.375
$ws.Range("L131").Value = 4483.3695
$ws.Range("M131").Value = 2138.625
$ws.Range("N131").Value = -14563.3695

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 56500
$ws.Range("J69").Value = 56500
$ws.Range("L69").Value = 56500
$ws.Range("N69").Value = -57998

$ws.Range("H72").Value = 56500
$ws.Range("J72").Value = 56500
$ws.Range("L72").Value = 169500
$ws.Range("N72").Value = -176988

$ws.Range("H102").Value = 10837.211
$ws.Range("I102").Value = 10837.211
$ws.Range("K102").Value = 10837.211
$ws.Range("M102").Value = -9215.210999999999

$ws.Range("H132").Value = 780165.4
$ws.Range("I132").Value = 4532
$ws.Range("J132").Value = 2137523.8
$ws.Range("K132").Value = 13596
$ws.Range("L132").Value = 6412571.399999999
$ws.Range("M132").Value = -11066
$ws.Range("N132").Value = -6417631.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3551
$ws.Range("I4").Value = 3551
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3551
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -3438
$ws.Range("N4").ClearContents() | Out-Null

$ws.Range("H7").Value = 6805.4375
$ws.Range("I7").Value = 6740.5835
$ws.Range("J7").Value = 7000
$ws.Range("K7").Value = 6740.5835
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = -6628.5835
$ws.Range("N7").Value = -7224

$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10452

$ws.Range("H28").Value = 3551
$ws.Range("I28").Value = 3551
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 3551
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -3319
$ws.Range("N28").ClearContents() | Out-Null

$ws.Range("H37").Value = 3551
$ws.Range("I37").Value = 3551
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 3551
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -3444
$ws.Range("N37").ClearContents() | Out-Null

$ws.Range("H40").Value = 2632.3635
$ws.Range("I40").Value = 2632.3635
$ws.Range("K40").Value = 2632.3635
$ws.Range("M40").Value = -2496.3635

$ws.Range("H93").Value = 8332.286
$ws.Range("I93").Value = 11036.2
$ws.Range("K93").Value = 11036.2
$ws.Range("M93").Value = -9788.200000000001

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents() | Out-Null

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents() | Out-Null

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents() | Out-Null
$ws.Range("N122").ClearContents() | Out-Null

$ws.Range("H126").Value = 6805.4375
$ws.Range("I126").Value = 6740.5835
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 20221.7505
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -17751.7505
$ws.Range("N126").Value = -25940

$ws.Range("H132").Value = 2796674
$ws.Range("I132").Value = 3786.7693
$ws.Range("K132").Value = 11360.3079
$ws.Range("M132").Value = -8830.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 527.75
$ws.Range("I100").Value = 527.75
$ws.Range("K100").Value = 1055.5
$ws.Range("M100").Value = -514.5

$ws.Range("H107").Value = 980.13336
$ws.Range("I107").Value = 1367.8889
$ws.Range("J107").Value = 398.5
$ws.Range("K107").Value = 4103.6667
$ws.Range("L107").Value = 1195.5
$ws.Range("M107").Value = -2183.6667
$ws.Range("N107").Value = -5035.5
